$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Eggs of Murex" row (row 44) entirely - shifts everything below up by one.
$ws.Rows.Item(44).Delete()

# After that deletion, the row that used to be "Shells NA" (originally row 51)
# is now at row 50. Delete it too.
$ws.Rows.Item(50).Delete()

# Set the W(kg) column (G) to 0 for all remaining species rows (35-50).
$ws.Range("G35:G50").Value = 0

# Update the Numb column (H) for Anadara transversa (row 38) from 1001 to -1.
$ws.Range("H38").Value = -1
